$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 858. This shifts the existing rows 858-899
# down to 859-900 (and Excel automatically extends the sheet dimension
# from A1:D899 to A1:D900).
$ws.Rows(858).Insert()

# Populate the newly inserted row with the new data point:
# 2026/02/22 (日), hour 8, ranking 201.
# Force the date column to stay plain text (matching the sheet's existing
# convention of storing dates as text rather than Excel date serials) by
# temporarily using a text number format, then copy the neighboring cell's
# style so the new cell ends up with the same (default) style as the rest
# of the column instead of a newly synthesized "text" style.
$ws.Range("A858").NumberFormat = "@"
$ws.Range("A858").Value = "2026/02/22"
$ws.Range("A858").Style = $ws.Range("A857").Style

$ws.Range("B858").Value = "日"
$ws.Range("C858").Value = 8
$ws.Range("D858").Value = 201
